$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 11).Value = -0.2336087822049224
$ws.Cells.Item(2, 10).Value = 0.4221894279166428
$ws.Cells.Item(2, 9).Value = -0.4043046464706727
$ws.Cells.Item(2, 8).Value = -0.4410326232298434
$ws.Cells.Item(2, 7).Value = -0.6804809672324722
$ws.Cells.Item(2, 6).Value = 0.0115444171491989
$ws.Cells.Item(2, 5).Value = -0.6733414736251095
$ws.Cells.Item(2, 4).Value = -0.2895456687149427
$ws.Cells.Item(2, 3).Value = -1.66950562919271
$ws.Cells.Item(2, 2).Value = -0.5351204465965399

$ws.Cells.Item(3, 11).Value = -0.191729189980311
$ws.Cells.Item(3, 10).Value = -0.5908441378320941
$ws.Cells.Item(3, 9).Value = 0.06730610189352677
$ws.Cells.Item(3, 8).Value = 0.723104312015092
$ws.Cells.Item(3, 7).Value = -0.1033897623722235
$ws.Cells.Item(3, 6).Value = -0.1401177391313942
$ws.Cells.Item(3, 5).Value = -0.3795660831340231
$ws.Cells.Item(3, 4).Value = 0.3124593012476481
$ws.Cells.Item(3, 3).Value = -0.3724265895266604
$ws.Cells.Item(3, 2).Value = 0.01136921538350649

$ws.Cells.Item(4, 11).Value = -0.04712738345727097
$ws.Cells.Item(4, 10).Value = 0.3870385516598326
$ws.Cells.Item(4, 9).Value = 0.44847509378141
$ws.Cells.Item(4, 8).Value = 0.0493601459296269
$ws.Cells.Item(4, 7).Value = 0.7075103856552477
$ws.Cells.Item(4, 6).Value = 1.363308595776813
$ws.Cells.Item(4, 5).Value = 0.5368145213894975
$ws.Cells.Item(4, 4).Value = 0.5000865446303268
$ws.Cells.Item(4, 3).Value = 0.2606382006276979
$ws.Cells.Item(4, 2).Value = 0.9526635850093691

$ws.Cells.Item(5, 11).Value = 0.6319705015114304
$ws.Cells.Item(5, 10).Value = 0.8590864059593566
$ws.Cells.Item(5, 9).Value = -0.1197149861197203
$ws.Cells.Item(5, 8).Value = 0.3144509489973832
$ws.Cells.Item(5, 7).Value = 0.3758874911189606
$ws.Cells.Item(5, 6).Value = -0.02322745673282245
$ws.Cells.Item(5, 5).Value = 0.6349227829927984
$ws.Cells.Item(5, 4).Value = 1.290720993114364
$ws.Cells.Item(5, 3).Value = 0.4642269187270481
$ws.Cells.Item(5, 2).Value = 0.4274989419678774

$ws.Cells.Item(6, 11).Value = 0.3702166863774111
$ws.Cells.Item(6, 10).Value = -0.05718027227819322
$ws.Cells.Item(6, 9).Value = 0.5762794571478953
$ws.Cells.Item(6, 8).Value = 0.8033953615958215
$ws.Cells.Item(6, 7).Value = -0.1754060304832554
$ws.Cells.Item(6, 6).Value = 0.2587599046338481
$ws.Cells.Item(6, 5).Value = 0.3201964467554255
$ws.Cells.Item(6, 4).Value = -0.07891850109635756
$ws.Cells.Item(6, 3).Value = 0.5792317386292632
$ws.Cells.Item(6, 2).Value = 1.235029948750828

$ws.Cells.Item(7, 11).Value = -0.02746841204387546
$ws.Cells.Item(7, 10).Value = -0.1501794284847013
$ws.Cells.Item(7, 9).Value = 0.01615752620344563
$ws.Cells.Item(7, 8).Value = -0.4112394324521587
$ws.Cells.Item(7, 7).Value = 0.2222202969739298
$ws.Cells.Item(7, 6).Value = 0.449336201421856
$ws.Cells.Item(7, 5).Value = -0.529465190657221
$ws.Cells.Item(7, 4).Value = -0.09529925554011737
$ws.Cells.Item(7, 3).Value = -0.03386271341853997
$ws.Cells.Item(7, 2).Value = -0.4329776612703231

$ws.Cells.Item(8, 11).Value = -0.4671716238107607
$ws.Cells.Item(8, 10).Value = -0.07055289228830908
$ws.Cells.Item(8, 9).Value = -0.1177866288434494
$ws.Cells.Item(8, 8).Value = -0.2404976452842752
$ws.Cells.Item(8, 7).Value = -0.07416069059612829
$ws.Cells.Item(8, 6).Value = -0.5015576492517326
$ws.Cells.Item(8, 5).Value = 0.1319020801743559
$ws.Cells.Item(8, 4).Value = 0.3590179846222821
$ws.Cells.Item(8, 3).Value = -0.6197834074567948
$ws.Cells.Item(8, 2).Value = -0.1856174723396913

$ws.Cells.Item(9, 11).Value = 0.07532753529099229
$ws.Cells.Item(9, 10).Value = -0.1447968545825803
$ws.Cells.Item(9, 9).Value = -0.3988488897063004
$ws.Cells.Item(9, 8).Value = -0.002230158183848807
$ws.Cells.Item(9, 7).Value = -0.0494638947389891
$ws.Cells.Item(9, 6).Value = -0.1721749111798149
$ws.Cells.Item(9, 5).Value = -0.005837956491668017
$ws.Cells.Item(9, 4).Value = -0.4332349151472724
$ws.Cells.Item(9, 3).Value = 0.2002248142788162
$ws.Cells.Item(9, 2).Value = 0.4273407187267424

$ws.Cells.Item(10, 11).Value = -0.4380952487963659
$ws.Cells.Item(10, 10).Value = -0.2591224913255812
$ws.Cells.Item(10, 9).Value = -0.02333404273891299
$ws.Cells.Item(10, 8).Value = -0.2434584326124856
$ws.Cells.Item(10, 7).Value = -0.4975104677362057
$ws.Cells.Item(10, 6).Value = -0.1008917362137541
$ws.Cells.Item(10, 5).Value = -0.1481254727688944
$ws.Cells.Item(10, 4).Value = -0.2708364892097202
$ws.Cells.Item(10, 3).Value = -0.1044995345215733
$ws.Cells.Item(10, 2).Value = -0.5318964931771777

$ws.Cells.Item(11, 11).Value = -0.5793653109721442
$ws.Cells.Item(11, 10).Value = -0.3716462008140141
$ws.Cells.Item(11, 9).Value = -0.3610856705546931
$ws.Cells.Item(11, 8).Value = -0.1821129130839084
$ws.Cells.Item(11, 7).Value = 0.0536755355027598
$ws.Cells.Item(11, 6).Value = -0.1664488543708128
$ws.Cells.Item(11, 5).Value = -0.4205008894945329
$ws.Cells.Item(11, 4).Value = -0.02388215797208129
$ws.Cells.Item(11, 3).Value = -0.07111589452722158
$ws.Cells.Item(11, 2).Value = -0.1938269109680474

$ws.Cells.Item(12, 11).Value = 0.1190458097769828
$ws.Cells.Item(12, 10).Value = -0.1906403594810787
$ws.Cells.Item(12, 9).Value = -0.4078492589559834
$ws.Cells.Item(12, 8).Value = -0.2001301487978533
$ws.Cells.Item(12, 7).Value = -0.1895696185385323
$ws.Cells.Item(12, 6).Value = -0.0105968610677476
$ws.Cells.Item(12, 5).Value = 0.2251915875189206
$ws.Cells.Item(12, 4).Value = 0.005067197645347965
$ws.Cells.Item(12, 3).Value = -0.2489848374783721
$ws.Cells.Item(12, 2).Value = 0.1476338940440795

$ws.Cells.Item(13, 11).Value = 0.5777029950204122
$ws.Cells.Item(13, 10).Value = -0.07272342619877098
$ws.Cells.Item(13, 9).Value = 0.08158070763323305
$ws.Cells.Item(13, 8).Value = -0.2281054616248284
$ws.Cells.Item(13, 7).Value = -0.4453143610997332
$ws.Cells.Item(13, 6).Value = -0.2375952509416031
$ws.Cells.Item(13, 5).Value = -0.2270347206822821
$ws.Cells.Item(13, 4).Value = -0.04806196321149736
$ws.Cells.Item(13, 3).Value = 0.1877264853751708
$ws.Cells.Item(13, 2).Value = -0.0323979044984018

$ws.Cells.Item(14, 11).Value = -0.0673936950407959
$ws.Cells.Item(14, 10).Value = 0.3396354339941604
$ws.Cells.Item(14, 9).Value = 0.5558004826123696
$ws.Cells.Item(14, 8).Value = -0.09462593860681362
$ws.Cells.Item(14, 7).Value = 0.05967819522519041
$ws.Cells.Item(14, 6).Value = -0.2500079740328711
$ws.Cells.Item(14, 5).Value = -0.4672168735077758
$ws.Cells.Item(14, 4).Value = -0.2594977633496457
$ws.Cells.Item(14, 3).Value = -0.2489372330903247
$ws.Cells.Item(14, 2).Value = -0.06996447561954

$ws.Cells.Item(15, 11).Value = 0.3609055008270807
$ws.Cells.Item(15, 10).Value = 0.5197544139825933
$ws.Cells.Item(15, 9).Value = 0.0207121968051735
$ws.Cells.Item(15, 8).Value = 0.4277413258401298
$ws.Cells.Item(15, 7).Value = 0.643906374458339
$ws.Cells.Item(15, 6).Value = -0.006520046760844223
$ws.Cells.Item(15, 5).Value = 0.1477840870711598
$ws.Cells.Item(15, 4).Value = -0.1619020821869017
$ws.Cells.Item(15, 3).Value = -0.3791109816618064
$ws.Cells.Item(15, 2).Value = -0.1713918715036764

$ws.Cells.Item(16, 11).Value = 2.760585277975261
$ws.Cells.Item(16, 10).Value = 0.5107824383638689
$ws.Cells.Item(16, 9).Value = 0.4546186173936249
$ws.Cells.Item(16, 8).Value = 0.6134675305491375
$ws.Cells.Item(16, 7).Value = 0.1144253133717177
$ws.Cells.Item(16, 6).Value = 0.521454442406674
$ws.Cells.Item(16, 5).Value = 0.7376194910248832
$ws.Cells.Item(16, 4).Value = 0.08719306980569996
$ws.Cells.Item(16, 3).Value = 0.241497203637704
$ws.Cells.Item(16, 2).Value = -0.06818896562035748

$ws.Cells.Item(17, 11).Value = 10.23793915510299
$ws.Cells.Item(17, 10).Value = 2.768788824245696
$ws.Cells.Item(17, 9).Value = 0.5189859846343043
$ws.Cells.Item(17, 8).Value = 0.4628221636640603
$ws.Cells.Item(17, 7).Value = 0.6216710768195729
$ws.Cells.Item(17, 6).Value = 0.1226288596421531
$ws.Cells.Item(17, 5).Value = 0.5296579886771094
$ws.Cells.Item(17, 4).Value = 0.7458230372953185
$ws.Cells.Item(17, 3).Value = 0.09539661607613537
$ws.Cells.Item(17, 2).Value = 0.2497007499081394

$ws.Cells.Item(18, 11).Value = -7.935912205685947
$ws.Cells.Item(18, 10).Value = 10.23534959445152
$ws.Cells.Item(18, 9).Value = 2.766199263594229
$ws.Cells.Item(18, 8).Value = 0.5163964239828361
$ws.Cells.Item(18, 7).Value = 0.4602326030125921
$ws.Cells.Item(18, 6).Value = 0.6190815161681047
$ws.Cells.Item(18, 5).Value = 0.1200392989906849
$ws.Cells.Item(18, 4).Value = 0.5270684280256412
$ws.Cells.Item(18, 3).Value = 0.7432334766438504
$ws.Cells.Item(18, 2).Value = 0.09280705542466716

$ws.Cells.Item(19, 11).Value = 0.2102926738762539
$ws.Cells.Item(19, 10).Value = -7.929378873102896
$ws.Cells.Item(19, 9).Value = 10.24188292703457
$ws.Cells.Item(19, 8).Value = 2.77273259617728
$ws.Cells.Item(19, 7).Value = 0.522929756565888
$ws.Cells.Item(19, 6).Value = 0.466765935595644
$ws.Cells.Item(19, 5).Value = 0.6256148487511566
$ws.Cells.Item(19, 4).Value = 0.1265726315737368
$ws.Cells.Item(19, 3).Value = 0.5336017606086931
$ws.Cells.Item(19, 2).Value = 0.7497668092269023

$ws.Cells.Item(20, 11).Value = 2.246646450696576
$ws.Cells.Item(20, 10).Value = 0.09907597889718328
$ws.Cells.Item(20, 9).Value = -8.040595568081965
$ws.Cells.Item(20, 8).Value = 10.1306662320555
$ws.Cells.Item(20, 7).Value = 2.661515901198209
$ws.Cells.Item(20, 6).Value = 0.4117130615868174
$ws.Cells.Item(20, 5).Value = 0.3555492406165733
$ws.Cells.Item(20, 4).Value = 0.5143981537720859
$ws.Cells.Item(20, 3).Value = 0.01535593659466611
$ws.Cells.Item(20, 2).Value = 0.4223850656296224

$ws.Cells.Item(21, 11).Value = -1.235129679813658
$ws.Cells.Item(21, 10).Value = 2.204900306709672
$ws.Cells.Item(21, 9).Value = 0.05732983491027921
$ws.Cells.Item(21, 8).Value = -8.08234171206887
$ws.Cells.Item(21, 7).Value = 10.0889200880686
$ws.Cells.Item(21, 6).Value = 2.619769757211305
$ws.Cells.Item(21, 5).Value = 0.3699669175999133
$ws.Cells.Item(21, 4).Value = 0.3138030966296693
$ws.Cells.Item(21, 3).Value = 0.4726520097851818
$ws.Cells.Item(21, 2).Value = -0.02639020739223796

$ws.Cells.Item(22, 11).Value = -1.270988795495144
$ws.Cells.Item(22, 10).Value = -1.209814617080943
$ws.Cells.Item(22, 9).Value = 2.230215369442386
$ws.Cells.Item(22, 8).Value = 0.0826448976429941
$ws.Cells.Item(22, 7).Value = -8.057026649336155
$ws.Cells.Item(22, 6).Value = 10.11423515080131
$ws.Cells.Item(22, 5).Value = 2.64508481994402
$ws.Cells.Item(22, 4).Value = 0.3952819803326282
$ws.Cells.Item(22, 3).Value = 0.3391181593623842
$ws.Cells.Item(22, 2).Value = 0.4979670725178967

$ws.Cells.Item(23, 11).Value = 0.7745058067040239
$ws.Cells.Item(23, 10).Value = -1.266950130452231
$ws.Cells.Item(23, 9).Value = -1.20577595203803
$ws.Cells.Item(23, 8).Value = 2.2342540344853
$ws.Cells.Item(23, 7).Value = 0.08668356268590799
$ws.Cells.Item(23, 6).Value = -8.05298798429324
$ws.Cells.Item(23, 5).Value = 10.11827381584423
$ws.Cells.Item(23, 4).Value = 2.649123484986935
$ws.Cells.Item(23, 3).Value = 0.3993206453755421
$ws.Cells.Item(23, 2).Value = 0.343156824405298

$ws.Cells.Item(24, 11).Value = 0.07303413297936051
$ws.Cells.Item(24, 10).Value = 0.6556286700130015
$ws.Cells.Item(24, 9).Value = -1.385827267143253
$ws.Cells.Item(24, 8).Value = -1.324653088729052
$ws.Cells.Item(24, 7).Value = 2.115376897794278
$ws.Cells.Item(24, 6).Value = -0.03219357400511441
$ws.Cells.Item(24, 5).Value = -8.171865120984263
$ws.Cells.Item(24, 4).Value = 9.999396679153206
$ws.Cells.Item(24, 3).Value = 2.530246348295912
$ws.Cells.Item(24, 2).Value = 0.2804435086845197

$ws.Cells.Item(25, 11).Value = 0.1277855351333463
$ws.Cells.Item(25, 10).Value = 0.03211049432398849
$ws.Cells.Item(25, 9).Value = 0.6147050313576294
$ws.Cells.Item(25, 8).Value = -1.426750905798625
$ws.Cells.Item(25, 7).Value = -1.365576727384424
$ws.Cells.Item(25, 6).Value = 2.074453259138906
$ws.Cells.Item(25, 5).Value = -0.07311721266048643
$ws.Cells.Item(25, 4).Value = -8.212788759639636
$ws.Cells.Item(25, 3).Value = 9.958473040497832
$ws.Cells.Item(25, 2).Value = 2.48932270964054

$ws.Cells.Item(26, 11).Value = -0.4475586702863481
$ws.Cells.Item(26, 10).Value = -0.1871396332876253
$ws.Cells.Item(26, 9).Value = -0.2828146740969831
$ws.Cells.Item(26, 8).Value = 0.2997798629366579
$ws.Cells.Item(26, 7).Value = -1.741676074219596
$ws.Cells.Item(26, 6).Value = -1.680501895805395
$ws.Cells.Item(26, 5).Value = 1.759528090717934
$ws.Cells.Item(26, 4).Value = -0.388042381081458
$ws.Cells.Item(26, 3).Value = -8.527713928060606
$ws.Cells.Item(26, 2).Value = 9.643547872076862

$ws.Cells.Item(27, 11).Value = -0.6684786343103865
$ws.Cells.Item(27, 10).Value = -1.503933630468878
$ws.Cells.Item(27, 9).Value = -1.243514593470155
$ws.Cells.Item(27, 8).Value = -1.339189634279513
$ws.Cells.Item(27, 7).Value = -0.7565950972458717
$ws.Cells.Item(27, 6).Value = -2.798051034402126
$ws.Cells.Item(27, 5).Value = -2.736876855987925
$ws.Cells.Item(27, 4).Value = 0.7031531305354048
$ws.Cells.Item(27, 3).Value = -1.444417341263988
$ws.Cells.Item(27, 2).Value = -9.584088888243137

$ws.Cells.Item(28, 11).Value = 0.1020883817579226
$ws.Cells.Item(28, 10).Value = 0.2021677416605441
$ws.Cells.Item(28, 9).Value = -0.6332872544979471
$ws.Cells.Item(28, 8).Value = -0.3728682174992243
$ws.Cells.Item(28, 7).Value = -0.4685432583085821
$ws.Cells.Item(28, 6).Value = 0.1140512787250589
$ws.Cells.Item(28, 5).Value = -1.927404658431195
$ws.Cells.Item(28, 4).Value = -1.866230480016994
$ws.Cells.Item(28, 3).Value = 1.573799506506335
$ws.Cells.Item(28, 2).Value = -0.573770965293057

$ws.Cells.Item(29, 11).Value = -0.2623493406516572
$ws.Cells.Item(29, 10).Value = 0.09143716542803909
$ws.Cells.Item(29, 9).Value = 0.1915165253306606
$ws.Cells.Item(29, 8).Value = -0.6439384708278306
$ws.Cells.Item(29, 7).Value = -0.3835194338291078
$ws.Cells.Item(29, 6).Value = -0.4791944746384656
$ws.Cells.Item(29, 5).Value = 0.1034000623951754
$ws.Cells.Item(29, 4).Value = -1.938055874761079
$ws.Cells.Item(29, 3).Value = -1.876881696346878
$ws.Cells.Item(29, 2).Value = 1.563148290176452

$ws.Cells.Item(30, 11).Value = -0.2227864824353526
$ws.Cells.Item(30, 10).Value = -0.370963872867798
$ws.Cells.Item(30, 9).Value = -0.01717736678810172
$ws.Cells.Item(30, 8).Value = 0.08290199311451979
$ws.Cells.Item(30, 7).Value = -0.7525530030439714
$ws.Cells.Item(30, 6).Value = -0.4921339660452486
$ws.Cells.Item(30, 5).Value = -0.5878090068546065
$ws.Cells.Item(30, 4).Value = -0.005214469820965406
$ws.Cells.Item(30, 3).Value = -2.04667040697722
$ws.Cells.Item(30, 2).Value = -1.985496228563019

$ws.Cells.Item(31, 11).Value = 0.1427612275365414
$ws.Cells.Item(31, 10).Value = -0.01468476146761361
$ws.Cells.Item(31, 9).Value = -0.1628621519000589
$ws.Cells.Item(31, 8).Value = 0.1909243541796373
$ws.Cells.Item(31, 7).Value = 0.2910037140822588
$ws.Cells.Item(31, 6).Value = -0.5444512820762324
$ws.Cells.Item(31, 5).Value = -0.2840322450775096
$ws.Cells.Item(31, 4).Value = -0.3797072858868674
$ws.Cells.Item(31, 3).Value = 0.2028872511467736
$ws.Cells.Item(31, 2).Value = -1.838568686009481

$ws.Cells.Item(32, 11).Value = -0.241400058615729
$ws.Cells.Item(32, 10).Value = 0.2592593979103418
$ws.Cells.Item(32, 9).Value = 0.1018134089061868
$ws.Cells.Item(32, 8).Value = -0.04636398152625851
$ws.Cells.Item(32, 7).Value = 0.3074225245534377
$ws.Cells.Item(32, 6).Value = 0.4075018844560592
$ws.Cells.Item(32, 5).Value = -0.427953111702432
$ws.Cells.Item(32, 4).Value = -0.1675340747037092
$ws.Cells.Item(32, 3).Value = -0.263209115513067
$ws.Cells.Item(32, 2).Value = 0.319385421520574

$ws.Cells.Item(33, 11).Value = 0.08445123591687528
$ws.Cells.Item(33, 10).Value = -0.1959066446185939
$ws.Cells.Item(33, 9).Value = 0.3047528119074769
$ws.Cells.Item(33, 8).Value = 0.1473068229033219
$ws.Cells.Item(33, 7).Value = -0.0008705675291234075
$ws.Cells.Item(33, 6).Value = 0.3529159385505728
$ws.Cells.Item(33, 5).Value = 0.4529952984531944
$ws.Cells.Item(33, 4).Value = -0.3824596977052969
$ws.Cells.Item(33, 3).Value = -0.1220406607065741
$ws.Cells.Item(33, 2).Value = -0.2177157015159319

$ws.Cells.Item(34, 11).Value = -0.1050777432881008
$ws.Cells.Item(34, 10).Value = 0.06689711455691082
$ws.Cells.Item(34, 9).Value = -0.2134607659785583
$ws.Cells.Item(34, 8).Value = 0.2871986905475125
$ws.Cells.Item(34, 7).Value = 0.1297527015433575
$ws.Cells.Item(34, 6).Value = -0.01842468888908786
$ws.Cells.Item(34, 5).Value = 0.3353618171906084
$ws.Cells.Item(34, 4).Value = 0.4354411770932299
$ws.Cells.Item(34, 3).Value = -0.4000138190652613
$ws.Cells.Item(34, 2).Value = -0.1395947820665385

$ws.Cells.Item(35, 11).Value = 0.2888921154092369
$ws.Cells.Item(35, 10).Value = -0.01697042433709459
$ws.Cells.Item(35, 9).Value = 0.155004433507917
$ws.Cells.Item(35, 8).Value = -0.1253534470275521
$ws.Cells.Item(35, 7).Value = 0.3753060094985187
$ws.Cells.Item(35, 6).Value = 0.2178600204943637
$ws.Cells.Item(35, 5).Value = 0.06968263006191837
$ws.Cells.Item(35, 4).Value = 0.4234691361416146
$ws.Cells.Item(35, 3).Value = 0.5235484960442361
$ws.Cells.Item(35, 2).Value = -0.3119065001142551

$ws.Cells.Item(36, 11).Value = -0.1448632037902657
$ws.Cells.Item(36, 10).Value = 0.4674667488970205
$ws.Cells.Item(36, 9).Value = 0.161604209150689
$ws.Cells.Item(36, 8).Value = 0.3335790669957007
$ws.Cells.Item(36, 7).Value = 0.0532211864602315
$ws.Cells.Item(36, 6).Value = 0.5538806429863024
$ws.Cells.Item(36, 5).Value = 0.3964346539821473
$ws.Cells.Item(36, 4).Value = 0.248257263549702
$ws.Cells.Item(36, 3).Value = 0.6020437696293982
$ws.Cells.Item(36, 2).Value = 0.7021231295320197

$ws.Cells.Item(37, 11).Value = 1.455535409161496
$ws.Cells.Item(37, 10).Value = 0.7671640239623843
$ws.Cells.Item(37, 9).Value = 1.379493976649671
$ws.Cells.Item(37, 8).Value = 1.073631436903339
$ws.Cells.Item(37, 7).Value = 1.245606294748351
$ws.Cells.Item(37, 6).Value = 0.9652484142128814
$ws.Cells.Item(37, 5).Value = 1.465907870738952
$ws.Cells.Item(37, 4).Value = 1.308461881734797
$ws.Cells.Item(37, 3).Value = 1.160284491302352
$ws.Cells.Item(37, 2).Value = 1.514070997382048

$ws.Cells.Item(38, 11).Value = 0.2348700177716323
$ws.Cells.Item(38, 10).Value = 0.511561173195739
$ws.Cells.Item(38, 9).Value = -0.1768102120033725
$ws.Cells.Item(38, 8).Value = 0.4355197406839137
$ws.Cells.Item(38, 7).Value = 0.1296572009375822
$ws.Cells.Item(38, 6).Value = 0.3016320587825939
$ws.Cells.Item(38, 5).Value = 0.02127417824712469
$ws.Cells.Item(38, 4).Value = 0.5219336347731955
$ws.Cells.Item(38, 3).Value = 0.3644876457690405
$ws.Cells.Item(38, 2).Value = 0.2163102553365951

$ws.Cells.Item(39, 10).Value = 0.2388379152847414
$ws.Cells.Item(39, 9).Value = 0.5155290707088481
$ws.Cells.Item(39, 8).Value = -0.1728423144902634
$ws.Cells.Item(39, 7).Value = 0.4394876381970228
$ws.Cells.Item(39, 6).Value = 0.1336250984506913
$ws.Cells.Item(39, 5).Value = 0.305599956295703
$ws.Cells.Item(39, 4).Value = 0.0252420757602338
$ws.Cells.Item(39, 3).Value = 0.5259015322863045
$ws.Cells.Item(39, 2).Value = 0.3684555432821496

$ws.Cells.Item(40, 9).Value = 0.3744780054549828
$ws.Cells.Item(40, 8).Value = 0.6511691608790895
$ws.Cells.Item(40, 7).Value = -0.03720222432002201
$ws.Cells.Item(40, 6).Value = 0.5751277283672642
$ws.Cells.Item(40, 5).Value = 0.2692651886209327
$ws.Cells.Item(40, 4).Value = 0.4412400464659443
$ws.Cells.Item(40, 3).Value = 0.1608821659304752
$ws.Cells.Item(40, 2).Value = 0.661541622456546

$ws.Cells.Item(41, 8).Value = 0.1336718235993181
$ws.Cells.Item(41, 7).Value = 0.4103629790234248
$ws.Cells.Item(41, 6).Value = -0.2780084061756867
$ws.Cells.Item(41, 5).Value = 0.3343215465115995
$ws.Cells.Item(41, 4).Value = 0.028459006765268
$ws.Cells.Item(41, 3).Value = 0.2004338646102796
$ws.Cells.Item(41, 2).Value = -0.07992401592518952

$ws.Cells.Item(42, 7).Value = 0.08834060834722172
$ws.Cells.Item(42, 6).Value = 0.3650317637713285
$ws.Cells.Item(42, 5).Value = -0.3233396214277831
$ws.Cells.Item(42, 4).Value = 0.2889903312595031
$ws.Cells.Item(42, 3).Value = -0.01687220848682837
$ws.Cells.Item(42, 2).Value = 0.1551026493581833

$ws.Cells.Item(43, 6).Value = 0.02147918641116785
$ws.Cells.Item(43, 5).Value = 0.2981703418352746
$ws.Cells.Item(43, 4).Value = -0.3902010433638369
$ws.Cells.Item(43, 3).Value = 0.2221289093234493
$ws.Cells.Item(43, 2).Value = -0.08373363042288225

$ws.Cells.Item(44, 5).Value = -0.00810701594554874
$ws.Cells.Item(44, 4).Value = 0.268584139478558
$ws.Cells.Item(44, 3).Value = -0.4197872457205535
$ws.Cells.Item(44, 2).Value = 0.1925427069667326

$ws.Cells.Item(45, 4).Value = -0.02625767267518964
$ws.Cells.Item(45, 3).Value = 0.2504334827489171
$ws.Cells.Item(45, 2).Value = -0.4379379024501944

$ws.Cells.Item(46, 3).Value = -0.04428949692388896
$ws.Cells.Item(46, 2).Value = 0.2324016585002178

$ws.Cells.Item(47, 2).Value = -0.09587373626955231

